# B3-and-B4-PowerPoint.pptx — theme1.xml <-> theme2.xml content swap.
#
# The authored change swaps the raw OOXML contents of ppt/theme/theme1.xml
# (which held the plain "Office Theme" / "Office" colour scheme) and
# ppt/theme/theme2.xml (which held the "Integral" / "Red Violet" colour
# scheme), while every relationship (slideMaster1, notesMaster1,
# presentation.xml) keeps pointing at the same target filenames.
#
# Net effect on the part that PowerPoint's object model exposes as "the"
# presentation theme (ppt/theme/theme2.xml, the one wired to SlideMaster /
# Presentation) is that its colour scheme changes from the "Red Violet"
# palette to the plain "Office" palette. fontScheme/fmtScheme are byte
# identical between the two original themes, so the colour swap is the
# only theme content that actually changes what renders.

function ConvertTo-ComColor([string]$hex) {
    # PowerPoint COM RGB() longs are 0x00BBGGRR (little endian of R,G,B).
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Office / "plain" theme palette — this is what theme1.xml already had,
# and what theme2.xml must end up with after the swap.
$officeColors = @(
    "000000", # 1  Dark1
    "FFFFFF", # 2  Light1
    "44546A", # 3  Dark2
    "E7E6E6", # 4  Light2
    "5B9BD5", # 5  Accent1
    "ED7D31", # 6  Accent2
    "A5A5A5", # 7  Accent3
    "FFC000", # 8  Accent4
    "4472C4", # 9  Accent5
    "70AD47", # 10 Accent6
    "0563C1", # 11 Hyperlink
    "954F72"  # 12 Followed Hyperlink
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Item($i).RGB = ConvertTo-ComColor $officeColors[$i - 1]
}
